$wb = $excel.ActiveWorkbook

# --- "main" sheet: bump the "count of keys" counter (27 -> 31) ---
$wsMain = $wb.Worksheets.Item('main')
$wsMain.Range("B2").Value = 31

# --- "keys" sheet: register the 4 new command keys ---
$wsKeys = $wb.Worksheets.Item('keys')
$wsKeys.Cells.Item(28, 1).Value = 'cmd_1'
$wsKeys.Cells.Item(28, 2).Value = 4
$wsKeys.Cells.Item(29, 1).Value = 'cmd_2'
$wsKeys.Cells.Item(29, 2).Value = 4
$wsKeys.Cells.Item(30, 1).Value = 'cmd_3'
$wsKeys.Cells.Item(30, 2).Value = 4
$wsKeys.Cells.Item(31, 1).Value = 'cmd_4'
$wsKeys.Cells.Item(31, 2).Value = 4

# --- "ru-RU" sheet: Russian translations for the new keys ---
$wsRu = $wb.Worksheets.Item('ru-RU')
$wsRu.Cells.Item(28, 1).Value = 'cmd_1'
$wsRu.Cells.Item(28, 2).Value = 'Ты открыл консоль?'
$wsRu.Cells.Item(29, 1).Value = 'cmd_2'
$wsRu.Cells.Item(29, 2).Value = 'Попробуй нажать в ней Alt+Enter'
$wsRu.Cells.Item(30, 1).Value = 'cmd_3'
$wsRu.Cells.Item(30, 2).Value = 'А потом ещё раз.'
$wsRu.Cells.Item(31, 1).Value = 'cmd_4'
$wsRu.Cells.Item(31, 2).Value = 'Если, конечно, не боишься.'

# --- "en-US" sheet: English translations for the new keys ---
$wsEn = $wb.Worksheets.Item('en-US')
$wsEn.Cells.Item(28, 1).Value = 'cmd_1'
$wsEn.Cells.Item(28, 2).Value = 'Do you open a console?'
$wsEn.Cells.Item(29, 1).Value = 'cmd_2'
$wsEn.Cells.Item(29, 2).Value = 'Try to press Alt+Enter in console window.'
$wsEn.Cells.Item(30, 1).Value = 'cmd_3'
$wsEn.Cells.Item(30, 2).Value = 'And then again.'
$wsEn.Cells.Item(31, 1).Value = 'cmd_4'
$wsEn.Cells.Item(31, 2).Value = 'Unless yor`re afraid, of course.'

# --- view state: selections land on the freshly-added rows ---
$wsKeys.Range('A28:B31').Select()
$wsEn.Range('B31').Select()

# ru-RU ends up as the active sheet/tab when the workbook was saved
$wsRu.Range('B32').Select()
$wsRu.Activate()
